# Auto-generated edit script applying cryptos price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.989.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.683.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -2.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.15%  "

$ws.Range("E10").Value = "  -2.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.57%  "

$ws.Range("E12").Value = "  -11.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.157.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.859.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.685.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("E19").Value = "  -5.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.96%  "

$ws.Range("E21").Value = "  -4.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.504"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0856"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("E36").Value = "  -5.31%  "

$ws.Range("E37").Value = "  -3.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "338.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.17"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.930"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "

$ws.Range("E42").Value = "  -5.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.616"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("E46").Value = "  -5.21%  "

$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  -3.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.088.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.46%  "
